$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.265.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3773"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07407"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8715"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.35"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.685"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.414"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.29"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07082"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008768"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.280.55"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.316"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.049.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.940"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.62"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.250"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.314"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08940"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7834"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.522"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.937"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01967"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05248"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.270"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.887"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.356"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +20.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1692"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.598"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5068"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.44"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.669"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06336"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.72%  "